$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 120
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value2 = $cell.Value2 * 100000
    }
}
